# The commit swaps the theme used by the deck's slide master (ppt/theme/theme1.xml,
# originally the "Integral" / "Red Violet" theme) with the theme used by the notes
# master (ppt/theme/theme2.xml, originally the built-in "Office Theme" / "Office"
# colour scheme). The font scheme and format scheme (fills/lines/effects) of the two
# themes are byte-for-byte identical, so the only observable difference is the
# 12-slot DrawingML colour scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink).
#
# We drive this through the slide's ThemeColorScheme, which maps 1:1 onto
# <a:clrScheme> of the master theme (ppt/theme/theme1.xml) that every slide in this
# deck inherits from (there is a single slide master). ThemeColor.RGB uses the
# standard COM BGR-packed integer (0x00BBGGRR), so each target hex colour below is
# converted accordingly.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Colors(1).RGB  = 0          # dk1      -> #000000
$tcs.Colors(2).RGB  = 16777215   # lt1      -> #FFFFFF
$tcs.Colors(3).RGB  = 6968388    # dk2      -> #44546A
$tcs.Colors(4).RGB  = 15132391   # lt2      -> #E7E6E6
$tcs.Colors(5).RGB  = 13998939   # accent1  -> #5B9BD5
$tcs.Colors(6).RGB  = 3243501    # accent2  -> #ED7D31
$tcs.Colors(7).RGB  = 10855845   # accent3  -> #A5A5A5
$tcs.Colors(8).RGB  = 49407      # accent4  -> #FFC000
$tcs.Colors(9).RGB  = 12874308   # accent5  -> #4472C4
$tcs.Colors(10).RGB = 4697456    # accent6  -> #70AD47
$tcs.Colors(11).RGB = 12673797   # hlink    -> #0563C1
$tcs.Colors(12).RGB = 7491477    # folHlink -> #954F72
